# edit.ps1 - apply the JAM/10036.xlsx changes described in the commit:
#   1. separate syst and stat errors in several experiments (already reflected
#      in the "%stat_u" / "%sys_u" shared strings; no further action needed
#      here - see the header row update below for the %dR_c -> %dR_u rename).
#   2. the R-model uncertainty (column U, "%dR_u") changes from 2% to 0.7%.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header U1: "%dR_c" -> "%dR_u" -----------------------------
$ws.Range("U1").Value = "%dR_u"

# --- 2. Update the R-model uncertainty values in column U (rows 2-20) ----
#     from 2 (%) down to 0.7 (%), per Christy instead of Simona's paper.
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 21).Value = 0.7
}

# --- 3. Update the view/selection state to match the saved workbook ------
# Scroll so column G is at the left edge and select U7, matching the
# recorded cursor position in the authored workbook.
$win = $excel.ActiveWindow
$win.ScrollColumn = 7
$win.ScrollRow = 1
$ws.Range("U7").Select() | Out-Null

# Restore the workbook tab-ratio (bookViews/workbookView@tabRatio) divider
# between the sheet tabs and the horizontal scroll bar.
$win.TabRatio = 0.993
